# Workbook contains a single data table on Sheet1. This weekly data refresh
# inserts 2 new daily records (rows) right above the existing row 152, which
# pushes every subsequent row down by 2 (the last 2 rows that fall off the
# bottom of the previous range reappear as brand-new rows 268 and 269). The
# 2 freshly inserted rows (new 152 and 153) are then populated with the new
# observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 152; this shifts rows 152:267
# down to 154:269 and grows the used range from A1:T267 to A1:T269.
$ws.Rows("152:153").Insert()

# Constant values shared by every data row in this table.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100104
$producto  = "Frutos de pepita"
$categoriaId = 100104005
$categoria = "Pera"
$unidad    = "$/caja 15 kilos empedrada"
$origen    = "Región de O'Higgins"
$kgUnidad  = 15

# --- New row 152: Packham's Triumph, Primera ---
$ws.Cells.Item(152, 1).Value = $mercadoId
$ws.Cells.Item(152, 2).Value = $mercado
$ws.Cells.Item(152, 3).Value = $region
$ws.Cells.Item(152, 4).Value = 44729
$ws.Cells.Item(152, 5).Value = $codreg
$ws.Cells.Item(152, 6).Value = $tipo
$ws.Cells.Item(152, 7).Value = $productoId
$ws.Cells.Item(152, 8).Value = $producto
$ws.Cells.Item(152, 9).Value = $categoriaId
$ws.Cells.Item(152, 10).Value = $categoria
$ws.Cells.Item(152, 11).Value = "Packham's Triumph"
$ws.Cells.Item(152, 12).Value = "Primera"
$ws.Cells.Item(152, 13).Value = 600
$ws.Cells.Item(152, 14).Value = 14000
$ws.Cells.Item(152, 15).Value = 15000
$ws.Cells.Item(152, 16).Value = 14500
$ws.Cells.Item(152, 17).Value = $unidad
$ws.Cells.Item(152, 18).Value = $origen
$ws.Cells.Item(152, 19).Value = 967
$ws.Cells.Item(152, 20).Value = $kgUnidad

# --- New row 153: Packham's Triumph, Segunda ---
$ws.Cells.Item(153, 1).Value = $mercadoId
$ws.Cells.Item(153, 2).Value = $mercado
$ws.Cells.Item(153, 3).Value = $region
$ws.Cells.Item(153, 4).Value = 44729
$ws.Cells.Item(153, 5).Value = $codreg
$ws.Cells.Item(153, 6).Value = $tipo
$ws.Cells.Item(153, 7).Value = $productoId
$ws.Cells.Item(153, 8).Value = $producto
$ws.Cells.Item(153, 9).Value = $categoriaId
$ws.Cells.Item(153, 10).Value = $categoria
$ws.Cells.Item(153, 11).Value = "Packham's Triumph"
$ws.Cells.Item(153, 12).Value = "Segunda"
$ws.Cells.Item(153, 13).Value = 300
$ws.Cells.Item(153, 14).Value = 12000
$ws.Cells.Item(153, 15).Value = 12000
$ws.Cells.Item(153, 16).Value = 12000
$ws.Cells.Item(153, 17).Value = $unidad
$ws.Cells.Item(153, 18).Value = $origen
$ws.Cells.Item(153, 19).Value = 800
$ws.Cells.Item(153, 20).Value = $kgUnidad
